# Weekly fruit/vegetable price update.
# A new weekly record is inserted at the top of the data block (row 23),
# pushing the existing rows 23-74 down to 24-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23; this shifts rows 23:74 down to 24:75
# and Excel carries the existing row formatting down with them.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with this week's record.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = "2022-08-31"
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 100112026
$ws.Range("G23").Value = "Haba"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 55
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = 14000
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Provincia de Limarí"
$ws.Range("P23").Value = 560
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
